$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 82, shifting existing rows 82:201 down to 83:202
$ws.Rows("82:82").Insert()

# Populate the newly inserted row 82 with a fresh data record
$ws.Range("A82").Value = 8
$ws.Range("B82").Value = "Terminal La Palmera de La Serena"
$ws.Range("C82").Value = "Coquimbo"
$ws.Range("D82").Value = 44557
$ws.Range("E82").Value = 4
$ws.Range("F82").Value = 100112012
$ws.Range("G82").Value = "Espinaca"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 2400
$ws.Range("K82").Value = 400
$ws.Range("L82").Value = 500
$ws.Range("M82").Value = 450
$ws.Range("N82").Value = "$/atado 300 a 500 gramos"
$ws.Range("O82").Value = "Provincia del Elquí"
$ws.Range("P82").Value = 900
$ws.Range("Q82").Value = 0.5
$ws.Range("R82").Value = "Hortaliza"
